# 2019 秋季招聘 — add resume/application-deadline column, drop the TODO sheet.

$wb = $excel.ActiveWorkbook

# --- Work on the main sheet first (while TODO still exists, indices are stable) ---
$ws = $wb.Worksheets.Item("秋招进度")

# Insert a new column B for "网申截止时间" (application deadline), shifting the
# rest of the table (投递岗位 .. offer进度) one column to the right.
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "网申截止时间"
$ws.Columns("B").ColumnWidth = 14.5

# 阿里's application deadline goes in the new column, on its row (row 3).
# Force text formatting first so the "2018.8.17" string isn't parsed as a date,
# then drop the formatting again so the cell keeps the plain/default style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2018.8.17"
$ws.Range("B3").ClearFormats()

# Shrink the header row's font (18pt -> 14pt) and tighten the row height to match.
$ws.Rows(1).Font.Size = 14
$ws.Rows(1).RowHeight = 18

# --- Remove the TODO sheet entirely ---
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("TODO").Delete()

# --- Make the remaining sheet the active / selected one ---
$null = $ws.Select()
$null = $ws.Range("B2").Select()
